$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column number format onto the two brand-new rows (114,115)
# before writing values, so they inherit style "s=2" like every other A-column cell
# instead of getting a freshly minted style entry.
$ws.Range("A113").Copy() | Out-Null
$ws.Range("A114:A115").PasteSpecial(-4122) | Out-Null

$ws.Range("A90").Value2 = 44232
$ws.Range("B90").Value2 = 14
$ws.Range("C90").Value2 = 63
$ws.Range("D90").Value2 = 394.4896681277395

$ws.Range("A91").Value2 = 44233
$ws.Range("B91").Value2 = 14
$ws.Range("C91").Value2 = 69
$ws.Range("D91").Value2 = 432.0601127113338

$ws.Range("A92").Value2 = 44234
$ws.Range("B92").Value2 = 6
$ws.Range("C92").Value2 = 69
$ws.Range("D92").Value2 = 432.0601127113338

$ws.Range("A93").Value2 = 44235
$ws.Range("B93").Value2 = 10
$ws.Range("C93").Value2 = 59
$ws.Range("D93").Value2 = 369.44270507201

$ws.Range("A94").Value2 = 44236
$ws.Range("B94").Value2 = 12
$ws.Range("C94").Value2 = 53
$ws.Range("D94").Value2 = 331.8722604884157

$ws.Range("A95").Value2 = 44237
$ws.Range("B95").Value2 = 1
$ws.Range("C95").Value2 = 55
$ws.Range("D95").Value2 = 344.3957420162805

$ws.Range("A96").Value2 = 44238
$ws.Range("B96").Value2 = 2
$ws.Range("C96").Value2 = 64
$ws.Range("D96").Value2 = 400.7514088916719

$ws.Range("A97").Value2 = 44239
$ws.Range("B97").Value2 = 8
$ws.Range("C97").Value2 = 75
$ws.Range("D97").Value2 = 469.6305572949279

$ws.Range("A98").Value2 = 44240
$ws.Range("B98").Value2 = 16
$ws.Range("C98").Value2 = 75
$ws.Range("D98").Value2 = 469.6305572949279

$ws.Range("A99").Value2 = 44241
$ws.Range("B99").Value2 = 15
$ws.Range("C99").Value2 = 77
$ws.Range("D99").Value2 = 482.1540388227927

$ws.Range("A100").Value2 = 44242
$ws.Range("B100").Value2 = 21
$ws.Range("C100").Value2 = 86
$ws.Range("D100").Value2 = 538.5097056981841

$ws.Range("A101").Value2 = 44243
$ws.Range("B101").Value2 = 12
$ws.Range("C101").Value2 = 98
$ws.Range("D101").Value2 = 613.6505948653726

$ws.Range("A102").Value2 = 44244
$ws.Range("B102").Value2 = 3
$ws.Range("C102").Value2 = 93
$ws.Range("D102").Value2 = 582.3418910457107

$ws.Range("A103").Value2 = 44245
$ws.Range("B103").Value2 = 11
$ws.Range("C103").Value2 = 86
$ws.Range("D103").Value2 = 538.5097056981841

$ws.Range("A104").Value2 = 44246
$ws.Range("B104").Value2 = 20
$ws.Range("C104").Value2 = 88
$ws.Range("D104").Value2 = 551.0331872260488

$ws.Range("A105").Value2 = 44247
$ws.Range("B105").Value2 = 11
$ws.Range("C105").Value2 = 99
$ws.Range("D105").Value2 = 619.912335629305

$ws.Range("A106").Value2 = 44248
$ws.Range("B106").Value2 = 8
$ws.Range("C106").Value2 = 97
$ws.Range("D106").Value2 = 607.3888541014402

$ws.Range("A107").Value2 = 44249
$ws.Range("B107").Value2 = 23
$ws.Range("C107").Value2 = 100
$ws.Range("D107").Value2 = 626.1740763932373

$ws.Range("A108").Value2 = 44250
$ws.Range("B108").Value2 = 23
$ws.Range("C108").Value2 = 105
$ws.Range("D108").Value2 = 657.4827802128992

$ws.Range("A109").Value2 = 44251
$ws.Range("B109").Value2 = 1
$ws.Range("C109").Value2 = 106
$ws.Range("D109").Value2 = 663.7445209768315

$ws.Range("A110").Value2 = 44252
$ws.Range("B110").Value2 = 14
$ws.Range("C110").Value2 = 117
$ws.Range("D110").Value2 = 732.6236693800877

$ws.Range("A111").Value2 = 44253
$ws.Range("B111").Value2 = 25
$ws.Range("C111").Value2 = 105
$ws.Range("D111").Value2 = 657.4827802128992

$ws.Range("A112").Value2 = 44254
$ws.Range("B112").Value2 = 12
$ws.Range("C112").Value2 = 97
$ws.Range("D112").Value2 = 607.3888541014402

$ws.Range("A113").Value2 = 44255
$ws.Range("B113").Value2 = 19
$ws.Range("C113").ClearContents() | Out-Null
$ws.Range("D113").ClearContents() | Out-Null

$ws.Range("A114").Value2 = 44256
$ws.Range("B114").Value2 = 11
$ws.Range("C114").ClearContents() | Out-Null
$ws.Range("D114").ClearContents() | Out-Null

$ws.Range("A115").Value2 = 44257
$ws.Range("B115").Value2 = 15
$ws.Range("C115").ClearContents() | Out-Null
$ws.Range("D115").ClearContents() | Out-Null
